# Add data for 2021-10-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name and header title to reflect new "through" date
$ws.Name = "Through 2021-10-10"
$ws.Range("B1").Value = "October 2021 (through October 10)"

# Update existing totals that changed
$ws.Range("B2").Value = 5
$ws.Range("L3").Value = 6
$ws.Range("B4").Value = 5
$ws.Range("V6").Value = 1
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 4
$ws.Range("L9").Value = 1
$ws.Range("V10").Value = 1
$ws.Range("L11").Value = 1
$ws.Range("AF12").Value = 2
$ws.Range("AZ13").Value = 2
$ws.Range("BJ13").Value = 2
$ws.Range("B19").Value = 1
$ws.Range("V25").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("B37").Value = 3
$ws.Range("AF37").Value = 2
